$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced back to
# Text format first, otherwise Excel auto-converts them to numeric cells
# (the source data keeps these as literal text, e.g. "523.38").
$textForceCells = @("D5", "D6", "D20", "D22", "D29", "D33", "D35", "D37", "D38", "D40", "D42", "D43", "D45", "D50", "D51")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "57.594.64"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "3.102.21"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "523.38"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "140.79"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.101.39"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").Value = "3.636.67"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "57.637.51"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "3.101.26"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").Value = "12.80"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").Value = "336.75"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  +2.65%  "
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("D29").Value = "6.52"
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +2.19%  "
$ws.Range("D33").Value = "20.92"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("D35").Value = "156.37"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("E36").Value = "  +2.68%  "
$ws.Range("D37").Value = "6.11"
$ws.Range("E37").Value = "  +2.86%  "
$ws.Range("D38").Value = "27.21"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D40").Value = "0.0664"
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.142.71"
$ws.Range("E41").Value = "  +1.22%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "3.94"
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("D43").Value = "0.685"
$ws.Range("E43").Value = "  +4.60%  "
$ws.Range("E44").Value = "  +11.54%  "
$ws.Range("D45").Value = "36.80"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "2.304.04"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("E49").Value = "  +5.50%  "
$ws.Range("D50").Value = "20.73"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").Value = "6.00"
$ws.Range("E51").Value = "  +2.16%  "
